$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 3.230985683306322
$ws.Cells.Item(2, 3).Value = 10.29869402782916
$ws.Cells.Item(2, 4).Value = 0.8054896365839992
$ws.Cells.Item(2, 5).Value = 8.660232485948974
$ws.Cells.Item(2, 7).Value = 22.99540183366846

# Row 3
$ws.Cells.Item(3, 2).Value = 0.00002074986032285508
$ws.Cells.Item(3, 3).Value = 0.002777888934908601
$ws.Cells.Item(3, 4).Value = 0.1575252929769615
$ws.Cells.Item(3, 5).Value = 0.496779210170732
$ws.Cells.Item(3, 7).Value = 0.657103141942925

# Row 4
$ws.Cells.Item(4, 2).Value = 1.459612070389937
$ws.Cells.Item(4, 3).Value = 1.667794583268128
$ws.Cells.Item(4, 4).Value = 0.1575252929769615
$ws.Cells.Item(4, 5).Value = 0.496779210170732
$ws.Cells.Item(4, 7).Value = 3.781711156805759

# Row 5
$ws.Cells.Item(5, 2).Value = 0.3048080303191223
$ws.Cells.Item(5, 3).Value = 1.667794583268128
$ws.Cells.Item(5, 4).Value = 26.21740644021617
$ws.Cells.Item(5, 5).Value = 9353990175.932438
$ws.Cells.Item(5, 7).Value = 9353990204.122446
